$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.839.96'
$ws.Range("E2").Value = '  -5.10%  '

$ws.Range("D3").Value = '1.817.19'
$ws.Range("E3").Value = '  -4.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '280.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5062'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3529'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.40%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.18%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.33%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8536'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.32%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07854'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.34%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.826.86'
$ws.Range("E13").Value = '  +68.86%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.032'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.63%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.88%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.43%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008089'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.46%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = '25.913.29'
$ws.Range("E20").Value = '  -4.95%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.757'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.37%  '

$ws.Range("B22").Value = 'Cosmos'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.26%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.122'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.91%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.02'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.09%  '

$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.164'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.25%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.677'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.91%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.75%  '

$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '108.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.64%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.311'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.34%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.211'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.92%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08768'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.12%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04802'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.01%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7392'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.89%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.128'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.60%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.835'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.36%  '

$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.003'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.091'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.30%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.445'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.62%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5355'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.72%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01854'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.48%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9815'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.54%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '113.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.37%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.217'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.64%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.213'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -12.23%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4706'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.70%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.004'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1377'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.44%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.228'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.46%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.45%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05901'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.02%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.485'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.12%  '
